$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 2797.565817734744, 2819.355250465005)
    3 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 6.15379541431027)
    4 = @(0.1169995834814548, 1.626987699542094, 3.223369029078222, 2797.565817734744, 2802.533174046845)
    5 = @(0.0006075818656279264, 0.04103571897497393, 18.71679738969934, 13.86384647080068, 32.62228716134062)
    6 = @(3.272327238179451, 2919.202174992006, 189.6080260415259, 2459690191846.092, 2459690194958.174)
    7 = @(0.2881169905109251, 1.626987699542094, 18.71679738969934, 14773364.14517103, 14773384.77707311)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
